# Remove thread sleeps from profile testcases:
#  - E2 (TestCase_B1) result flips from SKIP to FAIL
#  - Old rows 85-87 (TestCase_B84/B85/B86) get reshuffled: the "FAIL" marker
#    that used to live as a trailing result on row 85/86 moves to TestCase_B1,
#    and rows 85-87 become plain (no Results value), with a new TestCase_B86
#    (record-view) row appended, plus two brand-new rows 88/89
#    (TestCase_B87/B88) that exercise clicking the article/patent title.
#  - Trailing filler rows (formerly 88-99, blank except for a stray style)
#    are removed entirely; the sheet now ends at row 89.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. TestCase_B1 result flips from SKIP to FAIL ---------------------
$ws.Cells.Item(2, 5).Value = "FAIL"

# --- 2. Delete the old trailing filler rows (88-99), bottom-up ---------
for ($r = 99; $r -ge 88; $r--) {
    $ws.Rows.Item($r).EntireRow.Delete()
}

# --- 3. Reuse the already-correct visual formatting of row 84 for the --
#        five rewritten/new rows 85-89 (same border/fill/font/wrap combo)
$ws.Range("A84:E84").Copy($ws.Range("A85:E85"))
$ws.Range("A84:E84").Copy($ws.Range("A86:E86"))
$ws.Range("A84:E84").Copy($ws.Range("A87:E87"))
$ws.Range("A84:E84").Copy($ws.Range("A88:E88"))
$ws.Range("A84:E84").Copy($ws.Range("A89:E89"))

# --- 4. Row 85: TestCase_B84 / OPQA-613 / ALL search-results fields ----
$ws.Cells.Item(85, 1).Value = "TestCase_B84"
$ws.Cells.Item(85, 2).Value = "OPQA-613"
$ws.Cells.Item(85, 3).Value = "Verify that following fields get displayed correctly for a patent in ALL search results page:`na)Title`nb)Inventors`nc)Assignees`nd)Patent number`ne)Publication date`nf)Times cited count`ng)Comments count"
$ws.Cells.Item(85, 4).Value = "Y"
$ws.Cells.Item(85, 5).Value = $null
$ws.Rows.Item(85).RowHeight = 120

# --- 5. Row 86: TestCase_B85 / OPQA-614 / PATENTS search-results fields -
$ws.Cells.Item(86, 1).Value = "TestCase_B85"
$ws.Cells.Item(86, 2).Value = "OPQA-614"
$ws.Cells.Item(86, 3).Value = "Verify that following fields get displayed correctly for a patent in PATENTS search results page:`na)Title`nb)Inventors`nc)Assignees`nd)Patent number`ne)Publication date`nf)Times cited count`ng)Comments count"
$ws.Cells.Item(86, 4).Value = "Y"
$ws.Cells.Item(86, 5).Value = $null
$ws.Rows.Item(86).RowHeight = 120

# --- 6. Row 87: TestCase_B86 / OPQA-562 / record-view fields ------------
$ws.Cells.Item(87, 1).Value = "TestCase_B86"
$ws.Cells.Item(87, 2).Value = "OPQA-562"
$ws.Cells.Item(87, 3).Value = "Verify that following fields get displayed correctly for an patent in record view page:`na)Title`nb)Inventors`nc)Assignees`nd)Publication Date`ne)Publication Number`nf)Times Cited count`ng)Cited patents count`nh)Cited Articles count`ng)Comments count`ni)Abstract`nj)IPC Codes`nk)DETAILS link"
$ws.Cells.Item(87, 4).Value = "Y"
$ws.Cells.Item(87, 5).Value = $null
$ws.Rows.Item(87).RowHeight = 195

# --- 7. Row 88 (new): TestCase_B87 / OPQA-567 / click article title ----
#        (default row height - leave RowHeight untouched)
$ws.Cells.Item(88, 1).Value = "TestCase_B87"
$ws.Cells.Item(88, 2).Value = "OPQA-567"
$ws.Cells.Item(88, 3).Value = "Verify that record view page of a patent gets displayed when user clicks on article title in ALL search results page"
$ws.Cells.Item(88, 4).Value = "Y"
$ws.Cells.Item(88, 5).Value = $null

# --- 8. Row 89 (new): TestCase_B88 / OPQA-573 / click patent title ------
#        (default row height - leave RowHeight untouched)
$ws.Cells.Item(89, 1).Value = "TestCase_B88"
$ws.Cells.Item(89, 2).Value = "OPQA-573"
$ws.Cells.Item(89, 3).Value = "Verify that record view page of a patent gets displayed when user clicks a patent title in PATENTS search results page"
$ws.Cells.Item(89, 4).Value = "Y"
$ws.Cells.Item(89, 5).Value = $null
